$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in A2 (161 -> 170)
$ws.Range("A2").Value = 170

# Remove row 3 entirely (it previously held leftover C3/D3 values)
$ws.Rows("3:3").Delete()

# Move the active selection to A2 (was E8)
$ws.Range("A2").Select()
